$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values are text that may look numeric (e.g. "0.9990", "44.72").
# Force text storage via NumberFormat "@" so Excel keeps the exact literal
# (trailing zeros, thousand-dot formatting, etc.), then clear the temporary
# number-format override so the cell style matches the original (General/default).
$dCells = @("D2","D3","D4","D5","D6","D7","D8","D9","D10","D11","D12","D13","D14","D15","D16","D17","D18","D19","D20","D21","D22","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $dCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value = "28.680.69"
$ws.Range("D3").Value = "1.813.60"
$ws.Range("D4").Value = "1.002"
$ws.Range("D5").Value = "328.67"
$ws.Range("D6").Value = "0.9990"
$ws.Range("D7").Value = "0.4405"
$ws.Range("D8").Value = "0.3800"
$ws.Range("D9").Value = "44.72"
$ws.Range("D10").Value = "0.07709"
$ws.Range("D11").Value = "1.151"
$ws.Range("D12").Value = "22.81"
$ws.Range("D13").Value = "1.001"
$ws.Range("D14").Value = "6.343"
$ws.Range("D15").Value = "7.597"
$ws.Range("D16").Value = "1.813.12"
$ws.Range("D17").Value = "0.00001097"
$ws.Range("D18").Value = "0.06749"
$ws.Range("D19").Value = "81.59"
$ws.Range("D20").Value = "0.9988"
$ws.Range("D21").Value = "17.79"
$ws.Range("D22").Value = "6.336"
$ws.Range("D23").Value = "28.682.45"
$ws.Range("D24").Value = "11.89"
$ws.Range("D25").Value = "2.442"
$ws.Range("D26").Value = "20.80"
$ws.Range("D27").Value = "2.391"
$ws.Range("D28").Value = "152.86"
$ws.Range("D29").Value = "2.019.82"
$ws.Range("D30").Value = "1.286"
$ws.Range("D31").Value = "133.33"
$ws.Range("D32").Value = "3.970"
$ws.Range("D33").Value = "5.879"
$ws.Range("D34").Value = "0.09299"
$ws.Range("D35").Value = "0.2272"
$ws.Range("D36").Value = "12.32"
$ws.Range("D37").Value = "0.06400"
$ws.Range("D38").Value = "0.02349"
$ws.Range("D39").Value = "0.6697"
$ws.Range("D40").Value = "5.248"
$ws.Range("D41").Value = "1.212"
$ws.Range("D42").Value = "8.172"
$ws.Range("D43").Value = "1.449"
$ws.Range("D44").Value = "14.12"
$ws.Range("D45").Value = "0.9983"
$ws.Range("D46").Value = "0.6133"
$ws.Range("D47").Value = "3.812"
$ws.Range("D48").Value = "129.14"
$ws.Range("D49").Value = "2.059"
$ws.Range("D50").Value = "0.07117"
$ws.Range("D51").Value = "1.158"

foreach ($addr in $dCells) { $ws.Range($addr).ClearFormats() }

# Remaining changed cells (coin name, link, volume %) are plain text already.
$ws.Range("E2").Value = "  +2.18%  "
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("E5").Value = "  -2.81%  "
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("E7").Value = "  +3.23%  "
$ws.Range("E8").Value = "  +8.38%  "
$ws.Range("E9").Value = "  -2.19%  "
$ws.Range("E10").Value = "  +3.52%  "
$ws.Range("E11").Value = "  +0.18%  "
$ws.Range("E12").Value = "  -0.87%  "
$ws.Range("E13").Value = "  -0.03%  "
$ws.Range("E14").Value = "  +1.20%  "
$ws.Range("E15").Value = "  +4.36%  "
$ws.Range("E16").Value = "  -0.25%  "
$ws.Range("E17").Value = "  +1.15%  "
$ws.Range("E18").Value = "  +1.10%  "
$ws.Range("E19").Value = "  -0.50%  "
$ws.Range("E20").Value = "  -0.22%  "
$ws.Range("E21").Value = "  +2.91%  "
$ws.Range("E22").Value = "  -1.95%  "
$ws.Range("E23").Value = "  +2.07%  "
$ws.Range("E24").Value = "  -0.96%  "
$ws.Range("E25").Value = "  +2.18%  "
$ws.Range("E26").Value = "  +0.23%  "
$ws.Range("E27").Value = "  -4.13%  "
$ws.Range("E28").Value = "  -2.09%  "
$ws.Range("E29").Value = "  -0.17%  "
$ws.Range("E30").Value = "  -1.38%  "
$ws.Range("E31").Value = "  +0.50%  "
$ws.Range("E32").Value = "  -2.15%  "
$ws.Range("E33").Value = "  -1.47%  "
$ws.Range("E34").Value = "  +1.10%  "
$ws.Range("E35").Value = "  +4.83%  "
$ws.Range("E36").Value = "  -0.46%  "
$ws.Range("E37").Value = "  +2.04%  "
$ws.Range("E38").Value = "  -0.54%  "
$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("E39").Value = "  -0.66%  "
$ws.Range("B40").Value = "InternetComputer(DFINITY)"
$ws.Range("C40").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("E41").Value = "  -0.49%  "
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("E43").Value = "  -3.01%  "
$ws.Range("E44").Value = "  +0.36%  "
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("E47").Value = "  -1.54%  "
$ws.Range("E48").Value = "  +0.28%  "
$ws.Range("E49").Value = "  +0.57%  "
$ws.Range("E50").Value = "  +0.14%  "
$ws.Range("E51").Value = "  -1.92%  "
